# Updates cryptos list with latest price/volume snapshot
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text must not be re-interpreted as a number (keeps the
# dotted price strings, e.g. "29.888.51" or "243.00", stored as text).
$textCells = @(
    'D2',
    'D3',
    'D4',
    'D5',
    'D6',
    'D7',
    'D8',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D14',
    'D15',
    'D16',
    'D17',
    'D18',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D25',
    'D26',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D32',
    'D33',
    'D34',
    'D35',
    'D36',
    'D38',
    'D39',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)

# Plain text updates (coin names, links, already-textual percentages).
$textUpdates = @{
    'E2' = '  +0.09%  '
    'E3' = '  -0.16%  '
    'E4' = '  +0.17%  '
    'E5' = '  -0.73%  '
    'E6' = '  +0.92%  '
    'E7' = '  +0.25%  '
    'E8' = '  +3.00%  '
    'E9' = '  +0.96%  '
    'E10' = '  +3.59%  '
    'E11' = '  +1.22%  '
    'E12' = '  +5.55%  '
    'E13' = '  +0.04%  '
    'E14' = '  +3.42%  '
    'E15' = '  +0.38%  '
    'E16' = '  +0.13%  '
    'E17' = '  +0.34%  '
    'E18' = '  +1.79%  '
    'E20' = '  +1.09%  '
    'E21' = '  +0.86%  '
    'B22' = 'Chainlink'
    'C22' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    'E22' = '  +20.59%  '
    'B23' = 'Dai'
    'C23' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    'E23' = '  +0.27%  '
    'E24' = '  +0.17%  '
    'E25' = '  +15.39%  '
    'E26' = '  +2.14%  '
    'E28' = '  -0.11%  '
    'E29' = '  +4.26%  '
    'E30' = '  +2.07%  '
    'E31' = '  +2.00%  '
    'E32' = '  +5.09%  '
    'E33' = '  +3.38%  '
    'E34' = '  +1.16%  '
    'E35' = '  +1.61%  '
    'E36' = '  +1.96%  '
    'E38' = '  -0.93%  '
    'E39' = '  +1.05%  '
    'E41' = '  +1.95%  '
    'E42' = '  +1.00%  '
    'E43' = '  -2.12%  '
    'E44' = '  +1.91%  '
    'E45' = '  +0.25%  '
    'E46' = '  +2.72%  '
    'B47' = 'Maker'
    'C47' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'E47' = '  +5.33%  '
    'B48' = 'Quant'
    'C48' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'E48' = '  +2.05%  '
    'E49' = '  +2.83%  '
    'E50' = '  -0.05%  '
    'E51' = '  +8.82%  '
}

# Numeric-looking price strings that must be force-written as text.
$priceUpdates = @{
    'D2' = '29.888.51'
    'D3' = '1.897.12'
    'D4' = '1.002'
    'D5' = '0.7952'
    'D6' = '243.00'
    'D7' = '1.003'
    'D8' = '0.3214'
    'D9' = '26.27'
    'D10' = '0.07110'
    'D11' = '0.08070'
    'D12' = '0.7749'
    'D13' = '1.903.15'
    'D14' = '5.332'
    'D15' = '92.57'
    'D16' = '29.910.66'
    'D17' = '13.90'
    'D18' = '5.935'
    'D19' = '245.18'
    'D20' = '0.000007768'
    'D21' = '2.182.57'
    'D22' = '8.230'
    'D23' = '1.003'
    'D24' = '1.002'
    'D25' = '0.1627'
    'D26' = '9.349'
    'D27' = '165.34'
    'D28' = '18.76'
    'D29' = '2.097'
    'D30' = '1.383'
    'D31' = '1.542'
    'D32' = '4.491'
    'D33' = '0.05693'
    'D34' = '4.102'
    'D35' = '1.271'
    'D36' = '0.7406'
    'D38' = '2.696'
    'D39' = '0.01935'
    'D40' = '2.780'
    'D41' = '0.4467'
    'D42' = '72.37'
    'D43' = '5.879'
    'D44' = '0.8483'
    'D45' = '1.003'
    'D46' = '1.895'
    'D47' = '1.031.17'
    'D48' = '102.39'
    'D49' = '9.953'
    'D50' = '7.519'
    'D51' = '3.014'
}

foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $priceUpdates[$cell]
    $ws.Range($cell).ClearFormats()
}

foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}
